# Update the Hefei expo-listing workbook:
#  - bump several "want-to-go" counts (column F) on rows 2,5,6,9,10
#  - insert a brand-new event row ("合肥·Look Look动漫嘉年华") at row 11,
#    which pushes the previous rows 11-12 down to rows 12-13 and renumbers
#    the running index shown in column A
# Applied identically to both the "展览" (Exhibitions) and "全部类型"
# (All types) sheets, which mirror each other in this workbook.

$wb = $excel.ActiveWorkbook

function Update-ExpoSheet($ws) {
    # --- bump existing "want to go" counters ---
    $ws.Range("F2").Value = 2974
    $ws.Range("F5").Value = 6731
    $ws.Range("F6").Value = 1717
    $ws.Range("F9").Value = 59
    $ws.Range("F10").Value = 120

    # --- insert a new row at 11, shifting old rows 11-12 down to 12-13 ---
    $ws.Rows.Item(11).Insert()

    # Column A carries a bold/bordered/centred running-index style (same
    # style as the header + every other row). Grab that formatting from the
    # row above before writing the new index value, so the new cell matches
    # its neighbours instead of picking up the Insert's own ad-hoc style.
    $ws.Range("A10").Copy()
    $ws.Range("A11").PasteSpecial(-4122)
    $ws.Range("A11").Value = 10

    # Column B holds plain text dates like "2024-05-01"; force text so
    # Excel doesn't reinterpret it as a real date serial, then drop back to
    # the sheet's default (unstyled) cell format to match the other rows.
    $ws.Range("B11").NumberFormat = "@"
    $ws.Range("B11").Value = "2024-05-01"
    $ws.Range("B11").Style = "Normal"

    $ws.Range("C11").Value = "合肥·Look Look动漫嘉年华"
    $ws.Range("D11").Value = "新站区东方大道288号 少荃体育中心"
    $ws.Range("E11").Value = "2024.05.01 10:00-05.01 17:30"
    $ws.Range("F11").Value = 0
    $ws.Range("G11").Value = 29.9
    $ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=82311"
    $ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202403/jbUNtkAQ1709619599897.png"

    # The shifted-down rows keep their own data (Insert preserved it), but
    # the running index in column A needs to advance by one to stay in
    # sequence with the newly inserted row.
    $ws.Range("A12").Value = 11
    $ws.Range("A13").Value = 12
}

Update-ExpoSheet $wb.Worksheets.Item("展览")
Update-ExpoSheet $wb.Worksheets.Item("全部类型")

Write-Output "Hefei expo sheets updated"
